$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old "blank separator" row (row 173).
# This shifts the old rows 173 (blank), 174 (Total), 175 (sum [h]) and
# 176 (sum [working weeks]) down to rows 175, 176, 177 and 178.
$ws.Rows.Item(173).Insert()
$ws.Rows.Item(173).Insert()

# The row-insert operation copies formatting from the row above into the
# new rows, which also creates stray formatted (but empty) cells in
# column G that should not exist for these two new rows. Remove them.
$ws.Range("G173").Clear()
$ws.Range("G174").Clear()

# New row 173: a completed time entry (2014-08-04, 10:00 - 11:15).
$ws.Range("A173").Value = 2014
$ws.Range("B173").Value = 8
$ws.Range("C173").Value = 4
$ws.Range("D173").Value = 0.41666666666666669
$ws.Range("E173").Value = 0.46875
$ws.Range("F173").Formula = "=(E173-D173)*24*60"
$ws.Range("G173").Formula = "=F173/60"

# New row 174: a new time entry that only has a start time filled in so
# far (end time / duration are still blank).
$ws.Range("A174").Value = 2014
$ws.Range("B174").Value = 8
$ws.Range("C174").Value = 4
$ws.Range("D174").Value = 0.55208333333333337

# The "sum [min]" total now needs to include the newly added row 173.
$ws.Range("F176").Formula = "=SUM(F2:F173)"

# Select the new active cell, matching the edited workbook's view state.
[void]$ws.Range("E174").Select()
